# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E20) is re-ordered from newest-first
# (2002, 2001, 1912, 1911, 1910) to oldest-first (1910, 1911, 1912, 2001, 2002).
# Because the rows themselves are not moved, the "Valor Mora" amounts in
# column F have to follow their original period label so the data keeps
# pointing at the right period: the value that used to sit next to "2002"
# (33125) now belongs to the row labelled "1910", and the value that used
# to sit next to "1910" (8833) now belongs to the row labelled "2002".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-label the period column in chronological order.
$ws.Range("E16").Value = "1910"
$ws.Range("E17").Value = "1911"
$ws.Range("E18").Value = "1912"
$ws.Range("E19").Value = "2001"
$ws.Range("E20").Value = "2002"

# Keep each "Valor Mora" amount attached to its original period label.
$ws.Range("F16").Value = 8833
$ws.Range("F17").Value = 33125
$ws.Range("F18").Value = 33125
$ws.Range("F19").Value = 33125
$ws.Range("F20").Value = 33125
